$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0, styled bold + thin box border + centered/top aligned
$ws.Range("B1").Value = 0
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1

# A2 = 0, same style as B1 (copy the formatting instead of re-deriving it
# property-by-property, which would otherwise leave an extra unused style
# record behind and bloat the style table)
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B2 = "disconnected_elements" label, no special style
$ws.Range("B2").Value = "disconnected_elements"
